$d = $word.ActiveDocument

$d.Content.Find.Execute("601×5=3005", $true, $false, $false, $false, $false, $true, 1, $false, "140×6=840", 2) | Out-Null
$d.Content.Find.Execute("756×8=6048", $true, $false, $false, $false, $false, $true, 1, $false, "539×4=2156", 2) | Out-Null
$d.Content.Find.Execute("693×4=2772", $true, $false, $false, $false, $false, $true, 1, $false, "201×3=603", 2) | Out-Null
$d.Content.Find.Execute("884×8=7072", $true, $false, $false, $false, $false, $true, 1, $false, "118×8=944", 2) | Out-Null
$d.Content.Find.Execute("958×3=2874", $true, $false, $false, $false, $false, $true, 1, $false, "478×3=1434", 2) | Out-Null
$d.Content.Find.Execute("158×7=1106", $true, $false, $false, $false, $false, $true, 1, $false, "444×7=3108", 2) | Out-Null
$d.Content.Find.Execute("830×3=2490", $true, $false, $false, $false, $false, $true, 1, $false, "283×9=2547", 2) | Out-Null
$d.Content.Find.Execute("686×9=6174", $true, $false, $false, $false, $false, $true, 1, $false, "959×8=7672", 2) | Out-Null
$d.Content.Find.Execute("602×3=1806", $true, $false, $false, $false, $false, $true, 1, $false, "565×6=3390", 2) | Out-Null
$d.Content.Find.Execute("738×9=6642", $true, $false, $false, $false, $false, $true, 1, $false, "861×3=2583", 2) | Out-Null
$d.Content.Find.Execute("187×9=1683", $true, $false, $false, $false, $false, $true, 1, $false, "299×5=1495", 2) | Out-Null
$d.Content.Find.Execute("378×9=3402", $true, $false, $false, $false, $false, $true, 1, $false, "421×3=1263", 2) | Out-Null
$d.Content.Find.Execute("265×3=795", $true, $false, $false, $false, $false, $true, 1, $false, "414×3=1242", 2) | Out-Null
$d.Content.Find.Execute("685×7=4795", $true, $false, $false, $false, $false, $true, 1, $false, "288×3=864", 2) | Out-Null
$d.Content.Find.Execute("388×7=2716", $true, $false, $false, $false, $false, $true, 1, $false, "988×2=1976", 2) | Out-Null
$d.Content.Find.Execute("810×9=7290", $true, $false, $false, $false, $false, $true, 1, $false, "577×6=3462", 2) | Out-Null
$d.Content.Find.Execute("885×5=4425", $true, $false, $false, $false, $false, $true, 1, $false, "723×2=1446", 2) | Out-Null
$d.Content.Find.Execute("190×6=1140", $true, $false, $false, $false, $false, $true, 1, $false, "266×2=532", 2) | Out-Null
$d.Content.Find.Execute("149×6=894", $true, $false, $false, $false, $false, $true, 1, $false, "220×2=440", 2) | Out-Null
$d.Content.Find.Execute("124×8=992", $true, $false, $false, $false, $false, $true, 1, $false, "828×6=4968", 2) | Out-Null
$d.Content.Find.Execute("150×3=450", $true, $false, $false, $false, $false, $true, 1, $false, "981×9=8829", 2) | Out-Null
$d.Content.Find.Execute("461×2=922", $true, $false, $false, $false, $false, $true, 1, $false, "610×4=2440", 2) | Out-Null
$d.Content.Find.Execute("833×4=3332", $true, $false, $false, $false, $false, $true, 1, $false, "411×8=3288", 2) | Out-Null
$d.Content.Find.Execute("389×7=2723", $true, $false, $false, $false, $false, $true, 1, $false, "530×9=4770", 2) | Out-Null
$d.Content.Find.Execute("610×7=4270", $true, $false, $false, $false, $false, $true, 1, $false, "109×3=327", 2) | Out-Null
